# Applies the "Savesheet para diferentes bimestres" change:
# the per-cell lists of values (one per bimester) are collapsed down to
# the value that corresponds to the currently selected bimester.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "MEC-3B-Tec. Soldagem"
$ws.Range("E2").Value  = "-"

$ws.Range("B3").Value  = "-"
$ws.Range("D3").Value  = "MEC-3B-Trat. Térmicos"
$ws.Range("E3").Value  = "-"

$ws.Range("B4").Value  = "-"
$ws.Range("D4").Value  = "MEC-3B-Trat. Térmicos"

$ws.Range("B6").Value  = "-"
$ws.Range("C6").Value  = "-"
$ws.Range("D6").Value  = "MEC-3B-Trat. Térmicos"

$ws.Range("B7").Value  = "MEC-3B-Tec. Soldagem"
$ws.Range("D7").Value  = "MEC-3B-Trat. Térmicos"

$ws.Range("B8").Value  = "MEC-3B-Tec. Soldagem"
$ws.Range("D8").Value  = "MEC-3B-Tec. Soldagem"

$ws.Range("B10").Value = "MEC-3A-Tec. Soldagem"
$ws.Range("D10").Value = "-"
$ws.Range("E10").Value = "MEC-3A-Tec. Soldagem"

$ws.Range("D11").Value = "-"
$ws.Range("E11").Value = "MEC-3A-Tec. Soldagem"

$ws.Range("C15").Value = "-"

$ws.Range("B16").Value = "MEC-3A-Tec. Soldagem"
